# Insert a new data row for "Feria Lagunitas de Puerto Montt - Membrillo"
# at row 72, pushing the existing rows 72:91 down to 73:92.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$newRow = 72
$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = (Get-Date -Year 2022 -Month 6 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100104
$ws.Cells.Item($newRow, 8).Value = "Frutos de pepita"
$ws.Cells.Item($newRow, 9).Value = 100104003
$ws.Cells.Item($newRow, 10).Value = "Membrillo"
$ws.Cells.Item($newRow, 11).Value = "Champion"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 300
$ws.Cells.Item($newRow, 14).Value = 13000
$ws.Cells.Item($newRow, 15).Value = 14000
$ws.Cells.Item($newRow, 16).Value = 13500
$ws.Cells.Item($newRow, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($newRow, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value = 750
$ws.Cells.Item($newRow, 20).Value = 18
